$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the table with a new "2022" column (N), mirroring the formatting
# already used by the adjacent "2021" column (M).
$xlPasteFormats = -4122

$ws.Range("M3").Copy()
$ws.Range("N3").PasteSpecial($xlPasteFormats)

$ws.Range("M4").Copy()
$ws.Range("N4").PasteSpecial($xlPasteFormats)
$ws.Range("N4").Value = 2022

$ws.Range("M5").Copy()
$ws.Range("N5").PasteSpecial($xlPasteFormats)
$ws.Range("N5").Value = 98.8

$ws.Range("M6").Copy()
$ws.Range("N6").PasteSpecial($xlPasteFormats)
$ws.Range("N6").Value = 98

$ws.Range("M7").Copy()
$ws.Range("N7").PasteSpecial($xlPasteFormats)
$ws.Range("N7").Value = 96.9

$ws.Range("O4").Select()
